# Insert a new "total_concentrations" worksheet right after
# "input_concentrations", populated with the totals table
# (components H / PO4 / Cu across three sample rows).

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("input_concentrations")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "total_concentrations"

# Header row
$newSheet.Range("A1").Value = "H"
$newSheet.Range("B1").Value = "PO4"
$newSheet.Range("C1").Value = "Cu"

# Data rows
$newSheet.Range("A2").Value = 0.01
$newSheet.Range("B2").Value = 0.01
$newSheet.Range("C2").Value = 0.01

$newSheet.Range("A3").Value = 0.02
$newSheet.Range("B3").Value = 0.01
$newSheet.Range("C3").Value = 0.01

$newSheet.Range("A4").Value = 0.03
$newSheet.Range("B4").Value = 0.01
$newSheet.Range("C4").Value = 0.01
